$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Price) and E (Volume 1h) for rows 2-51,
# and every row's G (Hora) column moves from "4" to "5".
$rows = @(
    @{Row=2;  D="286.95";    E="4.27%"}
    @{Row=3;  D="28.42";     E="4.47%"}
    @{Row=4;  D="4.929";     E="1.79%"}
    @{Row=5;  D="0.06555";   E="2.61%"}
    @{Row=6;  D="7.255";     E="4.40%"}
    @{Row=7;  D="1.347";     E="12.38%"}
    @{Row=8;  D="0.9102";    E="3.82%"}
    @{Row=9;  D="0.1569";    E="3.81%"}
    @{Row=10; D="0.06747";   E="31.52%"}
    @{Row=11; D="0.07644";   E="1.69%"}
    @{Row=12; D="0.02982";   E="0.58%"}
    @{Row=13; D="0.08970";   E="-0.06%"}
    @{Row=14; D="0.001588";  E="0.85%"}
    @{Row=15; D="0.0006546"; E="2.50%"}
    @{Row=16; D="0.006023";  E="-2.72%"}
    @{Row=17; D="3.471";     E="-0.09%"}
    @{Row=18; D="3.395";     E="2.59%"}
    @{Row=19; D="2.241";     E="-0.44%"}
    @{Row=20; D="0.3157";    E="0.68%"}
    @{Row=21; D="0.1350";    E="0.08%"}
    @{Row=22; D="3.982";     E="1.47%"}
    @{Row=23; D="0.04464";   E="1.12%"}
    @{Row=24; D="0.1520";    E="10.12%"}
    @{Row=25; D="0.001185";  E="0.59%"}
    @{Row=26; D="0.004338";  E="12.51%"}
    @{Row=27; D="0.08000";   E="-55.78%"}
    @{Row=28; D="0.0001180"; E="-1.78%"}
    @{Row=29; D="0.0001635"; E="-15.73%"}
    @{Row=30; D="--";        E="--%"}
    @{Row=31; D="--";        E="--%"}
    @{Row=32; D="--";        E="--%"}
    @{Row=33; D="--";        E="--%"}
    @{Row=34; D="--";        E="--%"}
    @{Row=35; D="--";        E="--%"}
    @{Row=36; D="--";        E="--%"}
    @{Row=37; D="--";        E="--%"}
    @{Row=38; D="--";        E="--%"}
    @{Row=39; D="--";        E="--%"}
    @{Row=40; D="0.04162";   E="0.74%"}
    @{Row=41; D="0.006707";  E="-1.56%"}
    @{Row=42; D="0.1413";    E="20.51%"}
    @{Row=43; D="0.002139";  E="1.78%"}
    @{Row=44; D="0.01241";   E="7.97%"}
    @{Row=45; D="0.00005539"; E="6.51%"}
    @{Row=46; D="1.562";     E="-6.93%"}
    @{Row=47; D="0.01850";   E="-0.21%"}
    @{Row=48; D="--";        E="--%"}
    @{Row=49; D="--";        E="--%"}
    @{Row=50; D="--";        E="--%"}
    @{Row=51; D="--";        E="--%"}
)

foreach ($item in $rows) {
    $r = $item.Row

    # Force text storage (these columns hold text like "286.95" / "4.27%" /
    # "--" rather than numbers/percentages) and then reset the style back to
    # Normal so we don't leave a stray NumberFormat-driven style index behind.
    $dCell = $ws.Range("D$r")
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.Style = "Normal"

    $eCell = $ws.Range("E$r")
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
    $eCell.Style = "Normal"

    $gCell = $ws.Range("G$r")
    $gCell.NumberFormat = "@"
    $gCell.Value = "5"
    $gCell.Style = "Normal"
}
